$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 128; existing rows 128:241 shift down to 129:242.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with its data (copy of the constant
# columns from the former row 128 -- now row 129 -- plus the new weekly values).
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C128").Value = "Los Lagos"
$ws.Range("D128").Value = 44586
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = 100114014
$ws.Range("G128").Value = "Betarraga"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 1000
$ws.Range("K128").Value = 900
$ws.Range("L128").Value = 1000
$ws.Range("M128").Value = 950
$ws.Range("N128").Value = "`$/paquete 5 unidades"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 190
$ws.Range("Q128").Value = 5
$ws.Range("R128").Value = "Hortaliza"

# Keep the date formatting for column D consistent with the rest of the column.
$ws.Range("D128").NumberFormat = $ws.Range("D129").NumberFormat
